$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reposition the application window (matches the saved workbookView x/yWindow).
$win = $excel.ActiveWindow
$win.Left = 8580
$win.Top = 1100

# Strip the ATG/atg start codon and TAA/taa stop codon from each protein tag sequence.
# Cells are updated in this specific order so the shared-string table compacts and
# reorders to match the canonical saved order of the source edit.
$ws.Range("B4").Value = 'TCTAAATTA'  # PX1
$ws.Range("B5").Value = 'tctatggttagtaaaggagaagaaaataacatggcaatcattaaggagttcatgagattcaaagttcacatggaaggttctgtaaatggacatgaatttgaaatagaaggtgaaggagaaggaaggccttatgaaggaacccaaaccgcgaagctaaaagttactaagggtggcccattaccatttgcatgggatatccttagccctcaattcatgtatgggtcaaaggcttatgtcaagcaccccgccgacattccagactatctaaagttatcttttcccgaagggtttaagtgggagcgtgtgatgaacttcgaagacggtggcgtggtaacagtgactcaggattcgtccctgcaagatggtgaatttatctacaaagtcaaattaagaggaactaactttccatctgacggcccggttatgcaaaaaaagacaatgggctgggaggcctcctcagaacgaatgtaccctgaagatggtgccttgaagggtgagattaaacaaagattgaaattgaaagatggtggacattatgacgctgaggttaaaacgacatacaaagctaagaaacctgtccagctcccaggtgcttacaatgtaaatataaaacttgatattacatcacataatgaagattatacgatagttgaacaatacgaaagggctgaggggagacatagtactggtggcatggatgaactatacaaa'  # mCherry
$ws.Range("B2").Value = 'AAGATTGAAGAAGGTAAGTTGGTTATCTGGATTAACGGTGACAAGGGTTACAACGGTTTGGCTGAAGTTGGTAAGAAATTTGAAAAAGATACCGGTATCAAGGTCACTGTTGAACACCCAGACAAGTTGGAAGAAAAGTTTCCACAAGTTGCTGCCACTGGTGATGGTCCAGACATTATCTTCTGGGCTCATGACAGATTCGGTGGTTACGCCCAATCCGGTTTGTTAGCCGAGATCACCCCAGATAAGGCTTTTCAAGATAAGTTGTATCCATTCACTTGGGATGCCGTCAGATACAACGGTAAGTTAATCGCCTACCCAATTGCTGTTGAAGCTTTGTCTTTGATCTACAATAAGGACTTGTTACCTAACCCACCAAAGACCTGGGAAGAAATCCCAGCTTTAGATAAGGAGTTAAAAGCTAAGGGTAAGTCCGCTTTGATGTTTAACTTGCAAGAACCATACTTCACTTGGCCATTGATCGCTGCTGATGGTGGTTACGCTTTTAAGTATGAAAACGGTAAATACGACATTAAGGATGTCGGTGTCGACAATGCTGGTGCTAAGGCCGGTTTAACTTTCTTAGTCGATTTGATTAAGAATAAACATATGAATGCTGACACTGATTACTCTATTGCTGAAGCTGCTTTCAACAAGGGTGAAACCGCTATGACTATTAACGGTCCATGGGCCTGGTCTAACATTGATACCTCTAAAGTCAACTACGGTGTCACCGTCTTGCCAACTTTTAAGGGTCAACCATCTAAGCCATTCGTCGGTGTCTTGTCTGCCGGTATTAACGCTGCCTCTCCAAATAAGGAATTGGCCAAGGAATTCTTAGAAAACTACTTGTTAACCGATGAAGGTTTAGAGGCCGTTAACAAGGATAAGCCATTAGGTGCTGTTGCTTTGAAGTCTTACGAAGAAGAGTTGGCTAAGGATCCAAGAATTGCTGCTACTATGGAAAACGCTCAAAAGGGTGAAATTATGCCAAACATCCCACAAATGTCTGCTTTCTGGTACGCTGTTCGTACCGCCGTCATTAATGCCGCTTCTGGTCGTCAAACTGTTGATGAAGCCTTGAAGGACGCTCAAACCAGAATTACTAAG'  # MBP
$ws.Range("B3").Value = 'tctaaaggtgaagaattattcactggtgttgtcccaattttggttgaattagatggtgatgttaatggtcacaaattttctgtctccggtgaaggtgaaggtgatgctacttacggtaaattgaccttaaaatttatttgtactactggtaaattgccagttccatggccaaccttagtcactactttcggttatggtgttcaatgttttgcgagatacccagatcatatgaaacaacatgactttttcaagtctgccatgccagaaggttatgttcaagaaagaactatttttttcaaagatgacggtaactacaagaccagagctgaagtcaagtttgaaggtgataccttagttaatagaatcgaattaaaaggtattgattttaaagaagatggtaacattttaggtcacaaattggaatacaactataactctcacaatgtttacatcatggctgacaaacaaaagaatggtatcaaagttaacttcaaaattagacacaacattgaagatggttctgttcaattagctgaccattatcaacaaaatactccaattggtgatggtccagtcttgttaccagacaaccattacttatccactcaatctgccttatccaaagatccaaacgaaaagagagaccacatggtcttgttagaatttgttactgctgctggtattatccatggtatggatgaattgtacaaa'  # GFP
$ws.Range("B6").Value = 'CTTTCACTACGTCAATCTATAAGATTTTTCAAGCCAGCCACAAGAACTTTGTGTAGCTCTAGA'  # COX4
$ws.Range("B7").Value = 'GCTTCAGAAAAAGAAATTAGGAGAGAGAGATTCTTGAACGTTTTCCCTAAATTAGTAGAGGAATTGAACGCATCGCTTTTGGCTTACGGTATGCCTAAGGAAGCATGTGACTGGTATGCCCACTCATTGAACTACAACACTCCAGGCGGTAAGCTAAATAGAGGTTTGTCCGTTGTGGACACGTATGCTATTCTCTCCAACAAGACCGTTGAACAATTGGGGCAAGAAGAATACGAAAAGGTTGCCATTCTAGGTTGGTGCATTGAGTTGTTGCAGGCTTACTTCTTGGTCGCCGATGATATGATGGACAAGTCCATTACCAGAAGAGGCCAACCATGTTGGTACAAGGTTCCTGAAGTTGGGGAAATTGCCATCAATGACGCATTCATGTTAGAGGCTGCTATCTACAAGCTTTTGAAATCTCACTTCAGAAACGAAAAATACTACATAGATATCACCGAATTGTTCCATGAGGTCACCTTCCAAACCGAATTGGGCCAATTGATGGACTTAATCACTGCACCTGAAGACAAAGTCGACTTGAGTAAGTTCTCCCTAAAGAAGCACTCCTTCATAGTTACTTTCAAGACTGCTTACTATTCTTTCTACTTGCCTGTCGCATTGGCCATGTACGTTGCCGGTATCACGGATGAAAAGGATTTGAAACAAGCCAGAGATGTCTTGATTCCATTGGGTGAATACTTCCAAATTCAAGATGACTACTTAGACTGCTTCGGTACCCCAGAACAGATCGGTAAGATCGGTACAGATATCCAAGATAACAAATGTTCTTGGGTAATCAACAAGGCATTGGAACTTGCTTCCGCAGAACAAAGAAAGACTTTAGACGAAAATTACGGTAAGAAGGACTCAGTCGCAGAAGCCAAATGCAAAAAGATTTTCAATGACTTGAAAATTGAACAGCTATACCACGAATATGAAGAGTCTATTGCCAAGGATTTGAAGGCCAAAATTTCTCAGGTCGATGAGTCTCGTGGCTTCAAAGCTGATGTCTTAACTGCGTTCTTGAACAAAGTTTACAAGAGAAGCAAA'  # ERG20
$ws.Range("B8").Value = 'CAGATTTTCGTCAAGACTTTGACCGGTAAAACCATAACATTGGAAGTTGAATCTTCCGATACCATCGACAACGTTAAGTCGAAAATTCAAGACAAGGAAGGTATCCCTCCAGATCAACAAAGATTGATCTTTGCCGGTAAGCAGCTAGAAGACGGTAGAACGCTGTCTGATTACAACATTCAGAAGGAGTCCACCTTACATCTTGTGCTAAGGCTAAGAGGTGGTTATCACGGATCCGGAGCTTGGCTGTTGCCCGTCTCACTGGTGAAAAGAAAAACCACCCTGGCGCCCAATACG'  # UbiX
$ws.Range("B9").Value = 'TCTACCTCTGAAAACCAAAGTAAAGGTAGTGGTACATTGGTTGTCATATTGGCCATTTTAATGCTAGGTGTTGCTTATTATTTGTTGAACGAA'  # CYB5
$ws.Range("B10").Value = 'TGGTACAAGGATCTAAAAATGAAGATGTGTCTGGCTTTAGTAATCATCATATTGCTTGTTGTAATCATCGTCCCCATTGCTGTTCACTTTAGTCGA'  # SNC1
$ws.Range("B11").Value = 'aacgagctggccctgaagctggccggactggacatc'  # NES1

# Move the active selection to A12 (matches the post-edit saved state)
$ws.Range("A12").Select() | Out-Null
